# Zeiterfassung.xlsx - add two new weekly time-tracking blocks
# (06.11.17 and 07.11.17), each a copy of the existing block layout,
# separated from the prior data by a new grey-filled spacer row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1: 06.11.17 (rows 73-79), built from the 05.11.17 block (65-71) ---
$source = $ws.Range("A65:D71")
$target1 = $ws.Range("A73:D79")
$source.Copy($target1)

$ws.Range("A73").Value = "06.11.17"
$ws.Range("B78").Value = "1"
$ws.Range("C78").Value = "0.5"
$ws.Range("D78").Value = "1"
$ws.Range("B79").Value = 1
$ws.Range("C79").Value = 1
$ws.Range("D79").Value = 1

# --- Block 2: 07.11.17 (rows 81-87), same source layout ---
$target2 = $ws.Range("A81:D87")
$source.Copy($target2)

$ws.Range("A81").Value = "07.11.17"
$ws.Range("C85").Value = "1"

# --- Spacer rows 72 and 80: new grey fill style (distinct from the
#     older theme-coloured spacer rows used earlier in the sheet) ---
$spacer1 = $ws.Range("A72:D72")
$spacer1.Interior.PatternColor = 0
$spacer1.Interior.Color = 10855845
$spacer1.Font.Color = 0

$spacer2 = $ws.Range("A80:D80")
$spacer2.Interior.PatternColor = 0
$spacer2.Interior.Color = 10855845
$spacer2.Font.Color = 0

# --- Scroll / selection state so the view ends up near the new rows ---
$null = $ws.Range("C86").Select()
$excel.ActiveWindow.ScrollRow = 60
$excel.ActiveWindow.ScrollColumn = 1
